# The underlying edit swaps the entire contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the "Office Theme" palette that used to live in
# theme1.xml moves into theme2.xml, and the "Integral" palette that used to
# live in theme2.xml moves into theme1.xml (font scheme / format scheme are
# already byte-identical between the two parts, so only each theme's name +
# 12 scheme colors actually change).
#
# theme2.xml is the presentation's single reachable/"active" theme (the one
# behind $p.SlideMaster.Theme / Designs.Item(1)); theme1.xml only backs the
# Notes Master and isn't exposed as a separate object through the PowerPoint
# object model. So we drive the swap the same way a user would from the UI:
# by recoloring the active theme's 12 scheme colors to the target ("Office
# Theme") palette via ThemeColorScheme, in MSO slot order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function RgbVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$tcs.Item(1).RGB  = RgbVal 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RgbVal 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RgbVal 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RgbVal 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RgbVal 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RgbVal 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RgbVal 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RgbVal 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RgbVal 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RgbVal 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RgbVal 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RgbVal 0x95 0x4F 0x72   # folHlink
